# The underlying source data for this "Artfynd" export was re-sorted
# upstream; a handful of rows in the 46-58 range were re-ordered while
# their per-observation fields (id, species, coordinates, time, comment,
# substrate, ...) travelled together with them. Net effect on this sheet
# is a pure re-shuffle of whole data rows:
#   row46 <-> row47
#   row50 -> row51 -> row52 -> row50   (3-way rotation)
#   row57 <-> row58
#
# We copy whole rows (A:AY) so every column - including ones that differ
# only for a couple of these rows (I/J "Antal"/"Enhet") - moves as a unit.
# Columns Y ("Startdatum") and AA ("Slutdatum") are skipped on purpose:
# they hold the same literal text ("2026-01-20") in every one of these
# rows, so there is nothing to change there, and round-tripping a
# date-look-alike string through .Value tends to make Excel re-interpret
# it as a real date serial, which we want to avoid touching.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowBlocks($row) {
    return @{
        Left  = $ws.Range("A$($row):X$($row)").Value()
        Mid   = $ws.Range("Z$($row):Z$($row)").Value()
        Right = $ws.Range("AB$($row):AY$($row)").Value()
    }
}

function Set-RowBlocks($row, $blocks) {
    $ws.Range("A$($row):X$($row)").Value   = $blocks.Left
    $ws.Range("Z$($row):Z$($row)").Value   = $blocks.Mid
    $ws.Range("AB$($row):AY$($row)").Value = $blocks.Right
}

# Snapshot every row involved *before* writing anything, since several
# of them feed each other (the 50/51/52 rotation).
$src46 = Get-RowBlocks 46
$src47 = Get-RowBlocks 47
$src50 = Get-RowBlocks 50
$src51 = Get-RowBlocks 51
$src52 = Get-RowBlocks 52
$src57 = Get-RowBlocks 57
$src58 = Get-RowBlocks 58

# row46 <-> row47
Set-RowBlocks 46 $src47
Set-RowBlocks 47 $src46

# row50 -> row51 -> row52 -> row50
Set-RowBlocks 50 $src52
Set-RowBlocks 51 $src50
Set-RowBlocks 52 $src51

# row57 <-> row58
Set-RowBlocks 57 $src58
Set-RowBlocks 58 $src57
